$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Fill previously-blank cells in row 14 with the literal text "nan"
$ws.Range("B14").Value = "nan"
$ws.Range("C14").Value = "nan"
$ws.Range("D14").Value = "nan"
$ws.Range("E14").Value = "nan"
$ws.Range("F14").Value = "nan"
$ws.Range("G14").Value = "nan"
$ws.Range("H14").Value = "nan"
$ws.Range("I14").Value = "nan"
$ws.Range("J14").Value = "nan"
$ws.Range("K14").Value = "nan"
$ws.Range("M14").Value = "nan"
$ws.Range("P14").Value = "nan"

# Add new row 15 with data for a new maintenance event
$ws.Range("A15").Value = "'24"
$ws.Range("L15").Value = "17\9\2024"
$ws.Range("N15").Value = "تم عمل صيانه نصف سنويه"
$ws.Range("O15").Value = "تيم العمل"

# Row 15's remaining columns stay blank, like their row-14 counterparts
# (empty inline-string cells rather than wholly absent ones) - use the
# quote-prefix trick to get an empty *text* cell, then strip the style
# it leaves behind so formatting matches the untouched columns.
$blankCols15 = "B15","C15","D15","E15","F15","G15","H15","I15","J15","K15","M15","P15"
foreach ($addr in $blankCols15) {
    $ws.Range($addr).Value = "'"
}

$ws.Range("A15").Style = "Normal"
foreach ($addr in $blankCols15) {
    $ws.Range($addr).Style = "Normal"
}
